# Inserts the "Item Architecture" weapon-system journal entry right after the
# "Player Controller" paragraph, matching the target commit's diff.
$d = $word.ActiveDocument

# Locate the (bold, centered) "Player Controller" heading paragraph -- match case
# so we don't hit the earlier lower-case "player controller" mention in the prose.
$findRng = $d.Content
$found = $findRng.Find.Execute("Player Controller", $true, $true)
if (-not $found) {
    throw "Could not find the 'Player Controller' paragraph to anchor the insert"
}

$headingPara = $findRng.Paragraphs(1)
$followingPara = $headingPara.Next()

# Collapse a range to the very start of the paragraph that currently follows the
# heading, so the new paragraphs get spliced in between without disturbing it.
$insertRng = $d.Range($followingPara.Range.Start, $followingPara.Range.Start)

$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Item Architecture:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>So at first I used an inventory system hard coded with if statements to check what item I am holding and is fine on a low scale game. However once I start increasing the size of my game, this will become much harder to operate with as the code will become overly complex for no reason</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">To avoid this I can use abstract classes and inheritance  to let each item to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>onctrol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> itself.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Item class – contains abstract method use to use item in hand</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Then we have gun class to inherit item – has access to use method – This class will be parent of all guns and manages ammo and recoil</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>SingleShotGun</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> – single firing gun/ </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Shotgun- Override</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Can easily add more weapons </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# The fragment's final paragraph is empty (no <w:r>), so InsertXML adds it as its
# own standalone paragraph instead of merging into $followingPara -- leaving the
# pre-existing paragraph completely untouched right after it.
$insertRng.InsertXML($newParagraphsXml)

Write-Output "Inserted Item Architecture section after Player Controller"
